# Insert a new weekly data row for "Feria Lagunitas de Puerto Montt" / Betarraga
# at row 429, pushing the existing rows 429:461 down to 430:462, and
# populate the new row with this week's data (matching the formatting of the
# surrounding rows, which Excel's native row-insert already takes care of).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 429 (existing data shifts down one row).
$ws.Rows.Item(429).Insert()

# Fill in the new row 429 with the new record's values.
$ws.Cells.Item(429, 1).Value = 4
$ws.Cells.Item(429, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(429, 3).Value = "Los Lagos"
$ws.Cells.Item(429, 4).Value = 45021
$ws.Cells.Item(429, 5).Value = 10
$ws.Cells.Item(429, 6).Value = 100114014
$ws.Cells.Item(429, 7).Value = "Betarraga"
$ws.Cells.Item(429, 8).Value = "Sin especificar"
$ws.Cells.Item(429, 9).Value = "Primera"
$ws.Cells.Item(429, 10).Value = 250
$ws.Cells.Item(429, 11).Value = 1100
$ws.Cells.Item(429, 12).Value = 1100
$ws.Cells.Item(429, 13).Value = 1100
$ws.Cells.Item(429, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(429, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(429, 16).Value = 220
$ws.Cells.Item(429, 17).Value = 5
$ws.Cells.Item(429, 18).Value = "Hortaliza"
